$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J - copy H1's format (bold/border/centered style)
# then set the text, so the new header cells share the same style as existing headers.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2..26: I = 1 (constant), J = same value as H
for ($r = 2; $r -le 26; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
